# Applies the "Updated cryptos list" price/volume refresh described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.413.55"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "1.842.55"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'239.46"

$ws.Range("D6").Value = "'0.6278"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -0.71%  "

$ws.Range("D9").Value = "'0.2897"
$ws.Range("E9").Value = "  -0.12%  "

$ws.Range("D10").Value = "'24.95"
$ws.Range("E10").Value = "  +1.88%  "

$ws.Range("D11").Value = "'0.07717"
$ws.Range("E11").Value = "  -0.20%  "

$ws.Range("D12").Value = "1.842.34"
$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").Value = "'4.974"
$ws.Range("E13").Value = "  -0.48%  "

$ws.Range("D14").Value = "'0.6766"
$ws.Range("E14").Value = "  -0.73%  "

$ws.Range("D15").Value = "'0.00001027"
$ws.Range("E15").Value = "  -2.88%  "

$ws.Range("D16").Value = "'81.86"
$ws.Range("E16").Value = "  -0.32%  "

$ws.Range("D17").Value = "'6.243"
$ws.Range("E17").Value = "  +0.93%  "

$ws.Range("D18").Value = "29.465.72"
$ws.Range("E18").Value = "  +0.16%  "

$ws.Range("D19").Value = "'233.12"
$ws.Range("E19").Value = "  +1.60%  "

$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").Value = "'7.326"
$ws.Range("E22").Value = "  -2.16%  "

$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").Value = "'158.20"
$ws.Range("E24").Value = "  -0.76%  "

$ws.Range("D25").Value = "'8.496"

$ws.Range("E26").Value = "  -1.50%  "

$ws.Range("E27").Value = "  -0.95%  "

$ws.Range("D28").Value = "'0.07124"
$ws.Range("E28").Value = "  +9.57%  "

$ws.Range("D29").Value = "'1.468"
$ws.Range("E29").Value = "  +3.38%  "

$ws.Range("D30").Value = "'1.484"
$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.040"
$ws.Range("E31").Value = "  -1.41%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.047"
$ws.Range("E32").Value = "  -1.21%  "

$ws.Range("E33").Value = "  -0.67%  "

$ws.Range("D34").Value = "'1.141"
$ws.Range("E34").Value = "  -0.14%  "

$ws.Range("D35").Value = "'0.6958"
$ws.Range("E35").Value = "  -0.42%  "

$ws.Range("D36").Value = "'2.576"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("D37").Value = "'6.965"
$ws.Range("E37").Value = "  +3.00%  "

$ws.Range("D38").Value = "'0.01841"
$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("D39").Value = "'2.818"
$ws.Range("E39").Value = "  -0.67%  "

$ws.Range("D40").Value = "1.236.49"
$ws.Range("E40").Value = "  -2.40%  "

$ws.Range("D41").Value = "'0.9574"
$ws.Range("E41").Value = "  +5.21%  "

$ws.Range("D42").Value = "'1.0000"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").Value = "2.005.43"
$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("D44").Value = "'100.98"
$ws.Range("E44").Value = "  -0.30%  "

$ws.Range("D45").Value = "'65.50"
$ws.Range("E45").Value = "  -1.23%  "

$ws.Range("E46").Value = "  +0.44%  "

$ws.Range("D47").Value = "'1.730"
$ws.Range("E47").Value = "  -0.66%  "

$ws.Range("D48").Value = "'6.968"
$ws.Range("E48").Value = "  -1.55%  "

$ws.Range("D49").Value = "'8.903"
$ws.Range("E49").Value = "  -1.62%  "

$ws.Range("E50").Value = "  -2.28%  "

$ws.Range("D51").Value = "'0.3903"
$ws.Range("E51").Value = "  -1.76%  "
